$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("book")
$ws.Range("J2:J137").Formula = "PARENT"

$wsLib = $wb.Worksheets.Item("librarian")
$wsLib.Range("G4").Formula = '$2a$10$..4RvEGzO5/TiayeVSm1lOTHqo456ZPPCrf7G7.eyPvndgkICr/tq'
$wsLib.Range("G2").Formula = '$2a$10$Gitv.jdJOSpID30NmPEqn.IwH5CztayH4HbRUjxGDJKNM3DWCwMmy'
$wsLib.Range("G3").Formula = '$2a$10$GWg.rygxrh6caI2PaZTFvOreHBGaglUmkx1tyhutkxPN0QC1G1oT2'
$wsLib.Range("G5").Formula = '$2a$10$vrf9vKF0tT3xBQlEFbVjje.4LZgyCtRfBaIjCf4KO0QDQuvOKEijC'
